$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 46073
$ws.Range("C3").Value = 46073
$ws.Range("A4").Value = 'A 47653-2024'
$ws.Range("B4").Value = 45588
$ws.Range("C4").Value = 46073
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("R4").Value = 'Backklöver' + [char]13 + [char]10 + 'Nästrot' + [char]13 + [char]10 + 'Blåsippa'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/artfynd/A 47653-2024 artfynd.xlsx", "A 47653-2024")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/kartor/A 47653-2024 karta.png", "A 47653-2024")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomål/A 47653-2024 FSC-klagomål.docx", "A 47653-2024")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomålsmail/A 47653-2024 FSC-klagomål mail.docx", "A 47653-2024")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsyn/A 47653-2024 tillsynsbegäran.docx", "A 47653-2024")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsynsmail/A 47653-2024 tillsynsbegäran mail.docx", "A 47653-2024")'
$ws.Range("Z4").ClearContents()
$ws.Range("A5").Value = 'A 35197-2025'
$ws.Range("B5").Value = 45853
$ws.Range("C5").Value = 46073
$ws.Range("G5").Value = 0.9
$ws.Range("H5").Value = 1
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 0
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 0
$ws.Range("R5").Value = 'Grönsångare' + [char]13 + [char]10 + 'Tallticka' + [char]13 + [char]10 + 'Vintertagging'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/artfynd/A 35197-2025 artfynd.xlsx", "A 35197-2025")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/kartor/A 35197-2025 karta.png", "A 35197-2025")'
$ws.Range("U5").ClearContents()
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomål/A 35197-2025 FSC-klagomål.docx", "A 35197-2025")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomålsmail/A 35197-2025 FSC-klagomål mail.docx", "A 35197-2025")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsyn/A 35197-2025 tillsynsbegäran.docx", "A 35197-2025")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsynsmail/A 35197-2025 tillsynsbegäran mail.docx", "A 35197-2025")'
$ws.Range("Z5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/fåglar/A 35197-2025 prioriterade fågelarter.docx", "A 35197-2025")'
$ws.Range("A6").Value = 'A 37417-2023'
$ws.Range("B6").Value = 45155
$ws.Range("C6").Value = 46073
$ws.Range("G6").Value = 12.9
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 1
$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 1
$ws.Range("R6").Value = 'Knärot' + [char]13 + [char]10 + 'Ullticka' + [char]13 + [char]10 + 'Blåsippa'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/artfynd/A 37417-2023 artfynd.xlsx", "A 37417-2023")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/kartor/A 37417-2023 karta.png", "A 37417-2023")'
$ws.Range("U6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/knärot/A 37417-2023 karta knärot.png", "A 37417-2023")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomål/A 37417-2023 FSC-klagomål.docx", "A 37417-2023")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomålsmail/A 37417-2023 FSC-klagomål mail.docx", "A 37417-2023")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsyn/A 37417-2023 tillsynsbegäran.docx", "A 37417-2023")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsynsmail/A 37417-2023 tillsynsbegäran mail.docx", "A 37417-2023")'
$ws.Range("C7").Value = 46073
$ws.Range("C8").Value = 46073
$ws.Range("C9").Value = 46073
$ws.Range("C10").Value = 46073
$ws.Range("C11").Value = 46073
$ws.Range("A12").Value = 'A 55962-2023'
$ws.Range("B12").Value = 45240
$ws.Range("C12").Value = 46073
$ws.Range("G12").Value = 3.4
$ws.Range("A13").Value = 'A 23370-2025'
$ws.Range("B13").Value = 45791.70907407408
$ws.Range("C13").Value = 46073
$ws.Range("G13").Value = 3.8
$ws.Range("A14").Value = 'A 45406-2025'
$ws.Range("B14").Value = 45922.42936342592
$ws.Range("C14").Value = 46073
$ws.Range("G14").Value = 9.1
$ws.Range("A15").Value = 'A 35300-2025'
$ws.Range("B15").Value = 45854.41511574074
$ws.Range("C15").Value = 46073
$ws.Range("G15").Value = 2.2
$ws.Range("A16").Value = 'A 35198-2025'
$ws.Range("B16").Value = 45853
$ws.Range("C16").Value = 46073
$ws.Range("G16").Value = 1.2
$ws.Range("A17").Value = 'A 12156-2023'
$ws.Range("B17").Value = 44998.49157407408
$ws.Range("C17").Value = 46073
$ws.Range("G17").Value = 0.5
$ws.Range("A18").Value = 'A 32023-2023'
$ws.Range("B18").Value = 45119.49833333334
$ws.Range("C18").Value = 46073
$ws.Range("G18").Value = 3.1
$ws.Range("A19").Value = 'A 35193-2025'
$ws.Range("B19").Value = 45853
$ws.Range("C19").Value = 46073
$ws.Range("G19").Value = 1.9
$ws.Range("A20").Value = 'A 23250-2022'
$ws.Range("B20").Value = 44719
$ws.Range("C20").Value = 46073
$ws.Range("F20").Value = 'Naturvårdsverket'
$ws.Range("G20").Value = 1
$ws.Range("A21").Value = 'A 15732-2025'
$ws.Range("B21").Value = 45747
$ws.Range("C21").Value = 46073
$ws.Range("F21").Value = 'Kyrkan'
$ws.Range("G21").Value = 1.4
$ws.Range("A22").Value = 'A 8486-2026'
$ws.Range("B22").Value = 46064
$ws.Range("C22").Value = 46073
$ws.Range("G22").Value = 1
$ws.Range("C23").Value = 46073
$ws.Range("A24").Value = 'A 67005-2021'
$ws.Range("B24").Value = 44522
$ws.Range("C24").Value = 46073
$ws.Range("G24").Value = 1.3
$ws.Range("A25").Value = 'A 34202-2022'
$ws.Range("B25").Value = 44791.64837962963
$ws.Range("C25").Value = 46073
$ws.Range("G25").Value = 2
$ws.Range("A26").Value = 'A 4524-2024'
$ws.Range("B26").Value = 45327
$ws.Range("C26").Value = 46073
$ws.Range("G26").Value = 4.6
$ws.Range("A27").Value = 'A 65836-2021'
$ws.Range("B27").Value = 44517
$ws.Range("C27").Value = 46073
$ws.Range("G27").Value = 1.8
$ws.Range("C28").Value = 46073
$ws.Range("A29").Value = 'A 37415-2023'
$ws.Range("B29").Value = 45155
$ws.Range("C29").Value = 46073
$ws.Range("G29").Value = 6.6
$ws.Range("A30").Value = 'A 12146-2023'
$ws.Range("B30").Value = 44998.47842592592
$ws.Range("C30").Value = 46073
$ws.Range("F30").ClearContents()
$ws.Range("G30").Value = 3.1
$ws.Range("A31").Value = 'A 11989-2025'
$ws.Range("B31").Value = 45728.60074074074
$ws.Range("C31").Value = 46073
$ws.Range("F31").ClearContents()
$ws.Range("G31").Value = 9.6
